# Add 2022-Q3 data:
#  - insert a new row at the top of the "总计" (summary) sheet's data for
#    2022-Q3, shifting the existing quarters down by one row
#  - insert a new worksheet named "2022-Q3" right after "总计", built from
#    a copy of the old first quarterly sheet so formatting/styles match

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)      # "总计"
$wsQ2 = $wb.Worksheets.Item(2)     # currently "2022-Q2", used as a style template

# ---------------------------------------------------------------------
# 1) Update the "总计" summary sheet: shift existing rows down by one and
#    write the new 2022-Q3 summary row at the top.
# ---------------------------------------------------------------------

$summary = @(
    @("2022-Q3", 2, 0),
    @("2022-Q2", 4, 0.03),
    @("2021-Q4", 5, 0.14),
    @("2021-Q3", 2, 0.04),
    @("2021-Q2", 7, 0.37),
    @("2021-Q1", 4, 0.06),
    @("2020-Q4", 4, 0.08)
)

# row 8 is brand new (the sheet previously only went to row 7) - clone the
# "index" column formatting (bold + border, style used by A2:A7) onto it
# before filling in values, so it matches the rest of the column
$ws1.Range("A7").Copy()
$ws1.Range("A8").PasteSpecial(-4122)

for ($i = 0; $i -lt $summary.Count; $i++) {
    $r = $i + 2
    $row = $summary[$i]
    $ws1.Cells.Item($r, 1).Value = $i
    $ws1.Cells.Item($r, 2).Value = $row[0]
    $ws1.Cells.Item($r, 3).Formula = "$($row[1])"
    $ws1.Cells.Item($r, 4).Formula = "$($row[2])"
}

# ---------------------------------------------------------------------
# 2) Insert the new "2022-Q3" worksheet right after "总计" by duplicating
#    the old "2022-Q2" sheet (keeps header/style/border formatting) and
#    then overwriting its data with the 2022-Q3 fund holdings.
# ---------------------------------------------------------------------

$wsQ2.Copy($null, $ws1)
$newSheet = $wb.Worksheets.Item(2)
$newSheet.Name = "2022-Q3"

# the template sheet had 4 data rows (rows 2-5); the new sheet only needs 2
$newSheet.Rows.Item(5).Delete()
$newSheet.Rows.Item(4).Delete()

# force text formatting on the data cells B:G so numeric-looking strings
# (fund codes, percentages, ...) are kept as text, matching the source data
$textRange = $newSheet.Range("B2:G3")
$textRange.NumberFormat = "@"

$newSheet.Cells.Item(2, 1).Value = 0
$newSheet.Cells.Item(2, 2).Value = "012315"
$newSheet.Cells.Item(2, 3).Value = "创金合信港股通成长股票A"
$newSheet.Cells.Item(2, 4).Value = "0.08"
$newSheet.Cells.Item(2, 5).Value = "80.48"
$newSheet.Cells.Item(2, 6).Value = "3.28"
$newSheet.Cells.Item(2, 7).Value = "0.0026"
$newSheet.Cells.Item(2, 8).Value = 10

$newSheet.Cells.Item(3, 1).Value = 1
$newSheet.Cells.Item(3, 2).Value = "012316"
$newSheet.Cells.Item(3, 3).Value = "创金合信港股通成长股票C"
$newSheet.Cells.Item(3, 4).Value = "0.07"
$newSheet.Cells.Item(3, 5).Value = "80.48"
$newSheet.Cells.Item(3, 6).Value = "3.28"
$newSheet.Cells.Item(3, 7).Value = "0.0023"
$newSheet.Cells.Item(3, 8).Value = 10
